$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Paragraphs.Item(1).Range.Text = "2025-10-29 Wednesday"

# Update each math-problem cell in the table, row-major order
$t = $d.Tables.Item(1)
$values = @(
    "37+8=",
    "8+73=",
    "25-6=",
    "66-58=",
    "67+25=",
    "44+8=",
    "9+7=",
    "55+8=",
    "64-26=",
    "82-57=",
    "78+18=",
    "27+64=",
    "70-7=",
    "8+45=",
    "5+8=",
    "34-27=",
    "85-79=",
    "72-47=",
    "47+28=",
    "28+56=",
    "62-58=",
    "10-3=",
    "62-7=",
    "73-46=",
    "22+9=",
    "20-18=",
    "47+14=",
    "24+29=",
    "52-38=",
    "77+7=",
    "90-27=",
    "53-26=",
    "7+39=",
    "45+36=",
    "82-44=",
    "85-38=",
    "39+34=",
    "57-9=",
    "90-51=",
    "9+18=",
    "56+19=",
    "56-49=",
    "85-18=",
    "28+17=",
    "24+47=",
    "74-57=",
    "14+8=",
    "81-53=",
    "26+19=",
    "33+8=",
    "53-8=",
    "37+37=",
    "81-2=",
    "50-49=",
    "9+8=",
    "40-11=",
    "93-55=",
    "25+46=",
    "42-5=",
    "43+28=",
    "15+36=",
    "26+55=",
    "76+5=",
    "30-4=",
    "53+29=",
    "2+19=",
    "59+12=",
    "50-2=",
    "6+75=",
    "54-18=",
    "42-13=",
    "93-25=",
    "35+47=",
    "5+48=",
    "83-54=",
    "39+59=",
    "83-34=",
    "91-2=",
    "38+23=",
    "41-18=",
    "95-28=",
    "48+46=",
    "92-5=",
    "17+19=",
    "35-19=",
    "70-68=",
    "91-55=",
    "16+79=",
    "9+12=",
    "92-49=",
    "10-7=",
    "48+3=",
    "86+6=",
    "59+33=",
    "81-6=",
    "91-27=",
    "37+25=",
    "18+38=",
    "69+3=",
    "73-54="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count

if (($rows * $cols) -ne $values.Count) {
    Write-Output "WARNING: table has $rows x $cols = $($rows * $cols) cells but $($values.Count) values were supplied"
}

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Output "updated cells: $idx"